$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- 1. Insert a new bulleted paragraph "Method – get info" right before
#        the existing "toDoObject" bullet, carrying the bookmark that used
#        to live at the end of "valid date check" ------------------------
$idxToDo = Find-ParagraphIndex $d "toDoObject"
$pToDo = $d.Paragraphs.Item($idxToDo)

# Inserting before "toDoObject" clones its ListParagraph/numbered-list
# paragraph formatting onto a brand-new empty paragraph, without any of
# the proofErr spell-check markers that decorate "toDoObject" itself.
$pToDo.Range.InsertParagraphBefore()

$pNew = $d.Paragraphs.Item($idxToDo)

# Give the new paragraph its text. A throw-away trailing character ("Z")
# is appended so the bookmark we add next lands one character before the
# paragraph mark rather than exactly on it -- a collapsed bookmark placed
# precisely at the paragraph-end boundary gets mis-anchored back to the
# start of the document by this host, so we dodge that edge case and trim
# the placeholder off afterwards.
$pNew.Range.Text = "Method – get infoZ"
$pNewRange = $pNew.Range
$bmPos = $pNewRange.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($pNewRange.End - 2, $pNewRange.End - 1)
$placeholderRange.Text = ""

# --- 2. Re-join the two runs that spell out
#        "valid date check" + " (e.g. 31 of February)" into a single run.
#        (The bookmark that used to sit between them moved away in step 1,
#        since re-adding "_GoBack" elsewhere relocates the one bookmark of
#        that name instead of creating a duplicate.) ----------------------
$d.Content.Find.Execute(
    "valid date check (e.g. 31 of February)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "valid date check (e.g. 31 of February)", 2) | Out-Null
